# "Added values from this morning"
# Inserts two new columns (adapt.size.start / adapt.shape.start) into the
# MCMC settings table, backfills them for the existing scenario rows, and
# appends the new scenario rows gathered "this morning".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two blank columns before the old column N (accRate), pushing
#    accRate/ESS/Prop.fac./File names two columns to the right.
# ---------------------------------------------------------------------
$ws.Range("N1:O1").EntireColumn.Insert()
$ws.Range("N1:O1").ColumnWidth = 10.14

# ---------------------------------------------------------------------
# 2. New column headers (row 7, above the second data block) and the
#    per-row starting values for the already-present scenario rows.
# ---------------------------------------------------------------------
$ws.Range("O7").Value = "adapt.shape.start"
$ws.Range("N7").Value = "adapt.size.start"
$ws.Range("N7:O7").Font.Bold = $true

$ws.Range("N4:O4").NumberFormat = "0"

$sizeStartRows = 8,9,10,11,12,13,14,15
foreach ($r in $sizeStartRows) {
    $ws.Cells.Item($r, 14).Value = 100
    $ws.Cells.Item($r, 15).Value = 75
}

# ---------------------------------------------------------------------
# 3. Newly-added scenario rows (the "values from this morning").
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Belgium"
$ws.Range("B17").Value = "extloglin"
$ws.Range("C17").Value = -1.5
$ws.Range("D17").Value = -0.05
$ws.Range("E17").Value = -0.25
$ws.Range("F17").Value = -2.5
$ws.Range("G17").Value = -0.5
$ws.Range("H17").Value = -0.1
$ws.Range("I17").Value = 0.05
$ws.Range("J17").Value = -0.2
$ws.Range("K17").Value = 0.2
$ws.Range("L17").Value = 0.01
$ws.Range("M17").Value = 2000
$ws.Range("N17").Value = 100
$ws.Range("O17").Value = 75
$ws.Range("P17").NumberFormat = "0.0%"

$ws.Range("A18").Value = "Belgium"
$ws.Range("B18").Value = "extloglin"
$ws.Range("C18").Value = -1.5
$ws.Range("D18").Value = -0.05
$ws.Range("E18").Value = -0.25
$ws.Range("F18").Value = -2.5
$ws.Range("G18").Value = -1
$ws.Range("H18").Value = -0.1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = -0.15
$ws.Range("K18").Value = 0.15
$ws.Range("L18").Value = 0.01
$ws.Range("M18").Value = 2000
$ws.Range("N18").Value = 100
$ws.Range("O18").Value = 75
$ws.Range("P18").NumberFormat = "0.0%"

$ws.Range("A19").Value = "Belgium"
$ws.Range("B19").Value = "extloglin"
$ws.Range("C19").Value = -1.5
$ws.Range("D19").Value = -0.05
$ws.Range("E19").Value = -0.25
$ws.Range("F19").Value = -2.5
$ws.Range("G19").Value = -1
$ws.Range("H19").Value = -0.1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = -0.15
$ws.Range("K19").Value = 0.15
$ws.Range("L19").Value = 0.005
$ws.Range("M19").Value = 2000
$ws.Range("N19").Value = 200
$ws.Range("O19").Value = 100
$ws.Range("P19").Value = 0.134
$ws.Range("P19").NumberFormat = "0.0%"

$ws.Range("A21").Value = "Belgium"
$ws.Range("B21").Value = "extloglin"
$ws.Range("C21").Value = -1.5
$ws.Range("D21").Value = -0.05
$ws.Range("E21").Value = -0.25
$ws.Range("F21").Value = -2.5
$ws.Range("G21").Value = -1.5
$ws.Range("H21").Value = -0.1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = -0.1
$ws.Range("K21").Value = 0.1
$ws.Range("L21").Value = 0.01
$ws.Range("M21").Value = 2000
$ws.Range("N21").Value = 200
$ws.Range("O21").Value = 100
$ws.Range("P21").NumberFormat = "0.0%"

$ws.Range("A22").Value = "Belgium"
$ws.Range("B22").Value = "extloglin"
$ws.Range("C22").Value = -1.75
$ws.Range("D22").Value = -0.05
$ws.Range("E22").Value = -0.05
$ws.Range("F22").Value = -2.5
$ws.Range("G22").Value = -1
$ws.Range("H22").Value = -0.1
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = -0.1
$ws.Range("K22").Value = 0.1
$ws.Range("L22").Value = 0.01
$ws.Range("M22").Value = 2000
$ws.Range("N22").Value = 200
$ws.Range("O22").Value = 100
$ws.Range("P22").Value = 0.1945
$ws.Range("P22").NumberFormat = "0.0%"
$ws.Range("Q22").Value = 42.79082
$ws.Range("R22").Value = 57.03582
$ws.Range("S22").Value = 54.97821

$ws.Range("C23").Value = -1.75
$ws.Range("D23").Value = -0.05
$ws.Range("E23").Value = -0.05
$ws.Range("F23").Value = -2.5
$ws.Range("G23").Value = -1
$ws.Range("H23").Value = -0.1
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = -0.1
$ws.Range("K23").Value = 0.1
$ws.Range("L23").Value = 0.01
$ws.Range("M23").Value = 2000
$ws.Range("N23").Value = 100
$ws.Range("O23").Value = 75
$ws.Range("P23").Value = 0.219
$ws.Range("P23").NumberFormat = "0.0%"
$ws.Range("Q23").Value = 60.32554
$ws.Range("R23").Value = 60.73801
$ws.Range("S23").Value = 60.34873

$ws.Range("A24").Value = "Belgium"
$ws.Range("B24").Value = "extloglin"
$ws.Range("C24").Value = -1.75
$ws.Range("D24").Value = -0.05
$ws.Range("E24").Value = -0.05
$ws.Range("F24").Value = -2.5
$ws.Range("G24").Value = -1
$ws.Range("H24").Value = -0.1
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = -0.1
$ws.Range("K24").Value = 0.1
$ws.Range("L24").Value = 0.01
$ws.Range("M24").Value = 5000
$ws.Range("N24").Value = 100
$ws.Range("O24").Value = 75
$ws.Range("P24").Value = 0.1718
$ws.Range("P24").NumberFormat = "0%"
$ws.Range("Q24").Value = 119.6257
$ws.Range("R24").Value = 168.0576
$ws.Range("S24").Value = 128.7296

# ---------------------------------------------------------------------
# 4. Restore the selection to roughly where the author left off.
# ---------------------------------------------------------------------
$ws.Range("M27").Select()
